$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.842.97'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.19%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.857.26'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.77%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.31%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.78'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.90%  '

# Row 6
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.20%  '

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.59%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3656'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.58%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07188'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.31%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8894'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.49%  '

# Row 11
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.78%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07532'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.67%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.854.49'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.93%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '91.87'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.74%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.238'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.88%  '

# Row 16
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.36%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008541'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.17%  '

# Row 18
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.06%  '

# Row 19
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.19%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.881.71'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.23%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.014'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.56%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.094.50'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.45%  '

# Row 23
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.17%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.447'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.56%  '

# Row 25
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.10'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -3.70%  '

# Row 26
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.805'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.51%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.84'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.31%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.048'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -6.48%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '113.04'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.14%  '

# Row 30
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.38%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.668'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.96%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09224'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.96%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05094'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.76%  '

# Row 34
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7336'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.24%  '

# Row 35
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.149'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.65%  '

# Row 36
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.212'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +5.62%  '

# Row 37
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02008'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.76%  '

# Row 38
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.460'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.90%  '

# Row 39
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.074'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.87%  '

# Row 40
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5315'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.19%  '

# Row 41
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '118.07'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.49%  '

# Row 42
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.501'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.71%  '

# Row 43
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.389'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.17%  '

# Row 44
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1475'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.96%  '

# Row 45
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4634'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.45%  '

# Row 46
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9998'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.16%  '

# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.957'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.47%  '

# Row 48
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.558'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.19%  '

# Row 49
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.96'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.09%  '

# Row 50
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '62.87'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.49%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05932'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.43%  '
